$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.830.01'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '3.919.84'
$ws.Range("E3").Value = '  +4.16%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.44'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.74'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '3.917.90'
$ws.Range("E7").Value = '  +4.15%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  -1.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("E10").Value = '  -3.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.21'
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '4.571.94'
$ws.Range("E15").Value = '  +4.08%  '
$ws.Range("D16").Value = '3.941.57'
$ws.Range("E16").Value = '  +4.98%  '
$ws.Range("D17").Value = '68.946.97'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.49'
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.08'
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.15'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.37'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("E24").Value = '  +10.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.43'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.12'
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.13'
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").Value = '4.068.92'
$ws.Range("E31").Value = '  +4.11%  '
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '32.30'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D35").Value = '3.865.37'
$ws.Range("E35").Value = '  +4.42%  '
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("E37").Value = '  +3.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.141'
$ws.Range("E38").Value = '  +1.63%  '
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.321'
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '441.85'
$ws.Range("E42").Value = '  +3.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.00'
$ws.Range("E43").Value = '  -3.59%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.50'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.52'
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.847.31'
$ws.Range("E48").Value = '  +1.87%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.23'
$ws.Range("E49").Value = '  +11.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.84'
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0357'
$ws.Range("E51").Value = '  +1.14%  '
